$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Сотрудник", "Модель", "S\N (серийный номер)"),
    @("Молодцов Д.В.", "ASUS X4250LN-WX034H", "E4N0CX62583416A"),
    @("Родин Вадим", "Eee PC 1201NL", "A8OAAS275632"),
    @("Варибус П.И.", "R61e", "L3-ER768"),
    @("Варибус П.И.", "Comqaq 6715s", "CNU7450DTY8"),
    @("Лазарев Ю.П.", "ASUS K-43S", "EAN0CV006011418")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Column widths (compensate for the write-side +5/6-char padding applied when
# translating the COM "characters" ColumnWidth into the stored OOXML width, so
# the persisted <col width="..."> lands as close as possible to the source).
$ws.Columns.Item(1).ColumnWidth = 26.666666666666668
$ws.Columns.Item(2).ColumnWidth = 25.333333333333332
$ws.Columns.Item(3).ColumnWidth = 29.5

$ws.Range("B9").Select() | Out-Null
